$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E15").Value = 39.99
